$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Fzd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01948966666666667
$ws.Range("H2").Value = 0.058469
$ws.Range("I2").Value = 0.07096062449330311
$ws.Range("J2").Value = 0.07096062449330311
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.757644
$ws.Range("N2").Value = 2.272932
$ws.Range("O2").Value = 0.02401898721285653
$ws.Range("P2").Value = 0.02518910262217759
$ws.Range("Q2").Value = 0.014766229012
$ws.Range("R2").Value = 0.132896061108
$ws.Range("S2").Value = 0.001704402332320961
$ws.Range("T2").Value = 0.00178743445249562

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01948966666666667
$ws.Range("H3").Value = 0.058469
$ws.Range("I3").Value = 0.07096062449330311
$ws.Range("J3").Value = 0.07096062449330311
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 25.23919433333333
$ws.Range("N3").Value = 75.717583
$ws.Range("O3").Value = 0.8001381730141521
$ws.Range("P3").Value = 0.8391179183936208
$ws.Range("Q3").Value = 0.4919034844918889
$ws.Range("R3").Value = 4.427131360427
$ws.Range("S3").Value = 0.05677830443801484
$ws.Range("T3").Value = 0.05954433151273188

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Fzd1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01948966666666667
$ws.Range("H4").Value = 0.058469
$ws.Range("I4").Value = 0.07096062449330311
$ws.Range("J4").Value = 0.07096062449330311
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6151326666666667
$ws.Range("N4").Value = 1.845398
$ws.Range("O4").Value = 0.0195010633686494
$ws.Range("P4").Value = 0.02045108239083319
$ws.Range("Q4").Value = 0.01198873062911111
$ws.Range("R4").Value = 0.107898575662
$ws.Range("S4").Value = 0.001383807634922838
$ws.Range("T4").Value = 0.001451221578017518

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Fzd1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01948966666666667
$ws.Range("H5").Value = 0.058469
$ws.Range("I5").Value = 0.07096062449330311
$ws.Range("J5").Value = 0.07096062449330311
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5356743333333333
$ws.Range("N5").Value = 1.607023
$ws.Range("O5").Value = 0.01698205880675987
$ws.Range("P5").Value = 0.01780936132853939
$ws.Range("Q5").Value = 0.01044011419855555
$ws.Range("R5").Value = 0.093961027787
$ws.Range("S5").Value = 0.001205057498109678
$ws.Range("T5").Value = 0.001263763401700037

$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Fzd1"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01948966666666667
$ws.Range("H6").Value = 0.058469
$ws.Range("I6").Value = 0.07096062449330311
$ws.Range("J6").Value = 0.07096062449330311
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.395899500000001
$ws.Range("N6").Value = 8.791799000000001
$ws.Range("O6").Value = 0.139359717597582
$ws.Range("P6").Value = 0.09743253526482902
$ws.Range("Q6").Value = 0.08567461595516668
$ws.Range("R6").Value = 0.5140476957310001
$ws.Range("S6").Value = 0.009889052589934782
$ws.Range("T6").Value = 0.006913873548358045

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Fzd1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.255165
$ws.Range("H7").Value = 0.7654949999999999
$ws.Range("I7").Value = 0.9290393755066968
$ws.Range("J7").Value = 0.9290393755066969
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.757644
$ws.Range("N7").Value = 2.272932
$ws.Range("O7").Value = 0.02401898721285653
$ws.Range("P7").Value = 0.02518910262217759
$ws.Range("Q7").Value = 0.19332423126
$ws.Range("R7").Value = 1.73991808134
$ws.Range("S7").Value = 0.02231458488053557
$ws.Range("T7").Value = 0.02340166816968197

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt1"
$ws.Range("C8").Value = "Fzd1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.255165
$ws.Range("H8").Value = 0.7654949999999999
$ws.Range("I8").Value = 0.9290393755066968
$ws.Range("J8").Value = 0.9290393755066969
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 25.23919433333333
$ws.Range("N8").Value = 75.717583
$ws.Range("O8").Value = 0.8001381730141521
$ws.Range("P8").Value = 0.8391179183936208
$ws.Range("Q8").Value = 6.440159022064999
$ws.Range("R8").Value = 57.961431198585
$ws.Range("S8").Value = 0.7433598685761372
$ws.Range("T8").Value = 0.779573586880889

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt1"
$ws.Range("C9").Value = "Fzd1"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.255165
$ws.Range("H9").Value = 0.7654949999999999
$ws.Range("I9").Value = 0.9290393755066968
$ws.Range("J9").Value = 0.9290393755066969
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.6151326666666667
$ws.Range("N9").Value = 1.845398
$ws.Range("O9").Value = 0.0195010633686494
$ws.Range("P9").Value = 0.02045108239083319
$ws.Range("Q9").Value = 0.15696032689
$ws.Range("R9").Value = 1.41264294201
$ws.Range("S9").Value = 0.01811725573372656
$ws.Range("T9").Value = 0.01899986081281568

$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt1"
$ws.Range("C10").Value = "Fzd1"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.255165
$ws.Range("H10").Value = 0.7654949999999999
$ws.Range("I10").Value = 0.9290393755066968
$ws.Range("J10").Value = 0.9290393755066969
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5356743333333333
$ws.Range("N10").Value = 1.607023
$ws.Range("O10").Value = 0.01698205880675987
$ws.Range("P10").Value = 0.01780936132853939
$ws.Range("Q10").Value = 0.136685341265
$ws.Range("R10").Value = 1.230168071385
$ws.Range("S10").Value = 0.01577700130865019
$ws.Range("T10").Value = 0.01654559792683935

$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt1"
$ws.Range("C11").Value = "Fzd1"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.255165
$ws.Range("H11").Value = 0.7654949999999999
$ws.Range("I11").Value = 0.9290393755066968
$ws.Range("J11").Value = 0.9290393755066969
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.395899500000001
$ws.Range("N11").Value = 8.791799000000001
$ws.Range("O11").Value = 0.139359717597582
$ws.Range("P11").Value = 0.09743253526482902
$ws.Range("Q11").Value = 1.1216796959175
$ws.Range("R11").Value = 6.730078175505001
$ws.Range("S11").Value = 0.1294706650076472
$ws.Range("T11").Value = 0.09051866171647098
